$wb = $excel.ActiveWorkbook

# ---- Sheet: ALC ----
$ws = $wb.Worksheets.Item("ALC")
# Row 62
$ws.Range("H62").Value = 7018.273
$ws.Range("I62").Value = 7220.1
$ws.Range("K62").Value = 7220.1
$ws.Range("M62").Value = -6596.1
# Row 65
$ws.Range("H65").Value = 7018.273
$ws.Range("I65").Value = 7220.1
$ws.Range("K65").Value = 36100.5
$ws.Range("M65").Value = -32980.5
# Row 110
$ws.Range("H110").Value = 84696.664
$ws.Range("J110").Value = 84696.664
$ws.Range("L110").Value = 84696.664
$ws.Range("N110").Value = -92876.664
# Row 138
$ws.Range("H138").Value = 3046.1428
$ws.Range("I138").Value = 2443.5715
$ws.Range("J138").Value = 3166.6572
$ws.Range("K138").Value = 7330.7145
$ws.Range("L138").Value = 9499.971600000001
$ws.Range("M138").Value = -2190.7145
$ws.Range("N138").Value = -19779.9716

# ---- Sheet: ARM ----
$ws = $wb.Worksheets.Item("ARM")
# Row 2
$ws.Range("H2").Value = 544.9167
$ws.Range("J2").Value = 631.375
$ws.Range("L2").Value = 631.375
$ws.Range("N2").Value = -857.375
# Row 21
$ws.Range("H21").Value = 1161.2858
$ws.Range("I21").Value = 825.8
$ws.Range("J21").Value = 2000
$ws.Range("K21").Value = 825.8
$ws.Range("L21").Value = 2000
$ws.Range("M21").Value = -451.8
$ws.Range("N21").Value = -2748
# Row 32
$ws.Range("H32").Value = 287290.25
$ws.Range("I32").Value = 461974.97
$ws.Range("K32").Value = 461974.97
$ws.Range("M32").Value = -461687.97
# Row 45
$ws.Range("H45").Value = 2286.125
$ws.Range("I45").Value = 2597.5
$ws.Range("J45").Value = 1974.75
$ws.Range("K45").Value = 2597.5
$ws.Range("L45").Value = 1974.75
$ws.Range("M45").Value = -2220.5
$ws.Range("N45").Value = -2728.75
# Row 46
$ws.Range("H46").Value = 8633.666999999999
$ws.Range("J46").Value = 8633.666999999999
$ws.Range("L46").Value = 8633.666999999999
$ws.Range("N46").Value = -9271.666999999999
# Row 74
$ws.Range("H74").Value = 1857590.4
$ws.Range("I74").Value = 2139255.2
$ws.Range("K74").Value = 2139255.2
$ws.Range("M74").Value = -2138381.2
# Row 77
$ws.Range("H77").Value = 1857590.4
$ws.Range("I77").Value = 2139255.2
$ws.Range("K77").Value = 10696276
$ws.Range("M77").Value = -10691908
# Row 88
$ws.Range("H88").Value = 1681.1428
$ws.Range("J88").Value = 1811.3914
$ws.Range("L88").Value = 1811.3914
$ws.Range("N88").Value = -2623.3914
# Row 91
$ws.Range("H91").Value = 1681.1428
$ws.Range("J91").Value = 1811.3914
$ws.Range("L91").Value = 1811.3914
$ws.Range("N91").Value = -4619.3914
# Row 97
$ws.Range("H97").Value = 1891.5834
$ws.Range("I97").Value = 1891.5834
$ws.Range("J97").Value = 0
$ws.Range("K97").Value = 1891.5834
$ws.Range("L97").Value = 0
$ws.Range("N97").ClearContents()
# Row 113
$ws.Range("H113").Value = 75000
$ws.Range("J113").Value = 75000
$ws.Range("L113").Value = 75000
$ws.Range("N113").Value = -83678
# Row 116
$ws.Range("H116").Value = 544.9167
$ws.Range("J116").Value = 631.375
$ws.Range("L116").Value = 631.375
$ws.Range("N116").Value = -5219.375
# Row 122
$ws.Range("H122").Value = 2160.125
$ws.Range("I122").Value = 1756.2
$ws.Range("K122").Value = 5268.6
$ws.Range("M122").Value = -2818.6

# ---- Sheet: BSM ----
$ws = $wb.Worksheets.Item("BSM")
# Row 3
$ws.Range("H3").Value = 544.9167
$ws.Range("J3").Value = 631.375
$ws.Range("L3").Value = 631.375
$ws.Range("N3").Value = -859.375
# Row 86
$ws.Range("H86").Value = 1445.1702
$ws.Range("I86").Value = 1415.091
$ws.Range("J86").Value = 1516.0714
$ws.Range("K86").Value = 1415.091
$ws.Range("L86").Value = 1516.0714
$ws.Range("M86").Value = -292.0909999999999
$ws.Range("N86").Value = -3762.0714
# Row 89
$ws.Range("H89").Value = 1445.1702
$ws.Range("I89").Value = 1415.091
$ws.Range("J89").Value = 1516.0714
$ws.Range("K89").Value = 7075.455
$ws.Range("L89").Value = 7580.357
$ws.Range("M89").Value = -1459.455
$ws.Range("N89").Value = -18812.357
# Row 95
$ws.Range("H95").Value = 0
$ws.Range("J95").Value = 0
$ws.Range("N95").ClearContents()
# Row 107
$ws.Range("H107").Value = 887.9091
$ws.Range("I107").Value = 887.9091
$ws.Range("K107").Value = 887.9091
$ws.Range("M107").Value = 1032.0909
# Row 134
$ws.Range("H134").Value = 6671777.5
$ws.Range("I134").Value = 5424.7827
$ws.Range("J134").Value = 83334830
$ws.Range("K134").Value = 16274.3481
$ws.Range("L134").Value = 250004490
$ws.Range("M134").Value = -13739.3481
$ws.Range("N134").Value = -250009560

# ---- Sheet: CRP ----
$ws = $wb.Worksheets.Item("CRP")
# Row 116
$ws.Range("H116").Value = 78087
$ws.Range("I116").Value = 77674
$ws.Range("K116").Value = 77674
$ws.Range("M116").Value = -73085
# Row 134
$ws.Range("H134").Value = 2519.94
$ws.Range("I134").Value = 2273.7954
$ws.Range("J134").Value = 4325
$ws.Range("K134").Value = 6821.3862
$ws.Range("L134").Value = 12975
$ws.Range("M134").Value = -4286.3862
$ws.Range("N134").Value = -18045

# ---- Sheet: CUL ----
$ws = $wb.Worksheets.Item("CUL")
# Row 132
$ws.Range("H132").Value = 1416.4445
$ws.Range("I132").Value = 1049.6
$ws.Range("J132").Value = 1875
$ws.Range("K132").Value = 9446.4
$ws.Range("L132").Value = 16875
$ws.Range("M132").Value = -6916.4
$ws.Range("N132").Value = -21935

# ---- Sheet: GSM ----
$ws = $wb.Worksheets.Item("GSM")
# Row 5
$ws.Range("H5").Value = 204
$ws.Range("I5").Value = 204
$ws.Range("K5").Value = 204
$ws.Range("M5").Value = -92
# Row 70
$ws.Range("H70").Value = 57857.668
$ws.Range("I70").Value = 18087.5
$ws.Range("J70").Value = 82331.62
$ws.Range("K70").Value = 18087.5
$ws.Range("L70").Value = 82331.62
$ws.Range("M70").Value = -17817.5
$ws.Range("N70").Value = -82871.62
# Row 73
$ws.Range("H73").Value = 57857.668
$ws.Range("I73").Value = 18087.5
$ws.Range("J73").Value = 82331.62
$ws.Range("K73").Value = 18087.5
$ws.Range("L73").Value = 82331.62
$ws.Range("M73").Value = -17151.5
$ws.Range("N73").Value = -84203.62
# Row 80
$ws.Range("H80").Value = 3906.2812
$ws.Range("I80").Value = 2814.5
$ws.Range("J80").Value = 7181.625
$ws.Range("K80").Value = 2814.5
$ws.Range("L80").Value = 7181.625
$ws.Range("M80").Value = -1816.5
$ws.Range("N80").Value = -9177.625
# Row 83
$ws.Range("H83").Value = 3906.2812
$ws.Range("I83").Value = 2814.5
$ws.Range("J83").Value = 7181.625
$ws.Range("K83").Value = 14072.5
$ws.Range("L83").Value = 35908.125
$ws.Range("M83").Value = -9080.5
$ws.Range("N83").Value = -45892.125
# Row 95
$ws.Range("H95").Value = 36000
$ws.Range("J95").Value = 36000
$ws.Range("L95").Value = 36000
$ws.Range("N95").Value = -41492
# Row 97
$ws.Range("H97").Value = 1289.641
$ws.Range("I97").Value = 865.6
$ws.Range("K97").Value = 865.6
$ws.Range("M97").Value = -369.6
# Row 107
$ws.Range("H107").Value = 865.8333
$ws.Range("J107").Value = 1504.25
$ws.Range("L107").Value = 1504.25
$ws.Range("N107").Value = -5344.25
# Row 113
$ws.Range("H113").Value = 2486.9062
$ws.Range("I113").Value = 2368.5454
$ws.Range("J113").Value = 2747.3
$ws.Range("K113").Value = 2368.5454
$ws.Range("L113").Value = 2747.3
$ws.Range("M113").Value = -198.5454
$ws.Range("N113").Value = -7087.3

# ---- Sheet: LTW ----
$ws = $wb.Worksheets.Item("LTW")
# Row 5
$ws.Range("H5").Value = 12500
$ws.Range("J5").Value = 12500
$ws.Range("L5").Value = 12500
$ws.Range("N5").Value = -12726
# Row 22
$ws.Range("H22").Value = 1566
$ws.Range("I22").Value = 1297.3334
$ws.Range("J22").Value = 1834.6666
$ws.Range("K22").Value = 1297.3334
$ws.Range("L22").Value = 1834.6666
$ws.Range("M22").Value = -1002.3334
$ws.Range("N22").Value = -2424.6666
# Row 27
$ws.Range("H27").Value = 1566
$ws.Range("I27").Value = 1297.3334
$ws.Range("J27").Value = 1834.6666
$ws.Range("K27").Value = 1297.3334
$ws.Range("L27").Value = 1834.6666
$ws.Range("M27").Value = -1190.3334
$ws.Range("N27").Value = -2048.6666
# Row 61
$ws.Range("H61").Value = 6528.4443
$ws.Range("I61").Value = 4618.4116
$ws.Range("K61").Value = 4618.4116
$ws.Range("M61").Value = -4416.4116
# Row 95
$ws.Range("H95").Value = 57500
$ws.Range("J95").Value = 57500
$ws.Range("L95").Value = 57500
$ws.Range("N95").Value = -62992
# Row 100
$ws.Range("H100").Value = 4607
$ws.Range("I100").Value = 4607
$ws.Range("K100").Value = 4607
$ws.Range("M100").Value = -4066
# Row 113
$ws.Range("H113").Value = 6528.4443
$ws.Range("I113").Value = 4618.4116
$ws.Range("K113").Value = 4618.4116
$ws.Range("M113").Value = -2448.4116

# ---- Sheet: WVR ----
$ws = $wb.Worksheets.Item("WVR")
# Row 21
$ws.Range("H21").Value = 40000
$ws.Range("J21").Value = 40000
$ws.Range("L21").Value = 40000
$ws.Range("N21").Value = -40470
# Row 26
$ws.Range("H26").Value = 15200
$ws.Range("I26").Value = 10000
$ws.Range("J26").Value = 18666.666
$ws.Range("K26").Value = 10000
$ws.Range("L26").Value = 18666.666
$ws.Range("M26").Value = -9707
$ws.Range("N26").Value = -19252.666
# Row 30
$ws.Range("H30").Value = 35000
$ws.Range("J30").Value = 35000
$ws.Range("L30").Value = 35000
$ws.Range("N30").Value = -35214
# Row 35
$ws.Range("H35").Value = 40000
$ws.Range("J35").Value = 40000
$ws.Range("L35").Value = 40000
$ws.Range("N35").Value = -40580
# Row 95
$ws.Range("H95").Value = 0
$ws.Range("J95").Value = 0
$ws.Range("N95").ClearContents()
# Row 107
$ws.Range("H107").Value = 4127.852
$ws.Range("I107").Value = 3768.55
$ws.Range("K107").Value = 11305.65
$ws.Range("M107").Value = -9385.650000000001
# Row 126
$ws.Range("H126").Value = 2849.15
$ws.Range("I126").Value = 2740.353
$ws.Range("J126").Value = 3465.6667
$ws.Range("K126").Value = 8221.059000000001
$ws.Range("L126").Value = 10397.0001
$ws.Range("M126").Value = -5751.059000000001
$ws.Range("N126").Value = -15337.0001
